$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.432.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "'1.996.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.82%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'330.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("D8").Value = "'0.4185"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.20%  "
$ws.Range("D9").Value = "'53.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "'0.08836"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.54%  "
$ws.Range("D11").Value = "'1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.96%  "
$ws.Range("D12").Value = "'23.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.43%  "
$ws.Range("D13").Value = "'2.045.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "'8.021"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.48%  "
$ws.Range("D15").Value = "'6.478"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.49%  "
$ws.Range("D16").Value = "'96.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.64%  "
$ws.Range("D17").Value = "'1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'0.00001106"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.99%  "
$ws.Range("D19").Value = "'0.06625"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'19.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.56%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'5.963"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.21%  "
$ws.Range("D23").Value = "'29.485.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'11.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.288"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "'2.291.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'157.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'6.536"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'20.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.19%  "
$ws.Range("E30").Value = "  -7.68%  "
$ws.Range("D31").Value = "'126.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("D32").Value = "'1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.99%  "
$ws.Range("D33").Value = "'0.09920"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("D34").Value = "'1.555"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.67%  "
$ws.Range("D35").Value = "'5.836"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.78%  "
$ws.Range("D36").Value = "'3.771"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'9.587"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.16%  "
$ws.Range("D38").Value = "'0.02452"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.34%  "
$ws.Range("D39").Value = "'0.06379"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.12%  "
$ws.Range("D40").Value = "'1.286"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "'11.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.66%  "
$ws.Range("D42").Value = "'0.6498"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.85%  "
$ws.Range("D43").Value = "'0.2067"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.79%  "
$ws.Range("D44").Value = "'1.007"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'0.6317"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.90%  "
$ws.Range("D46").Value = "'13.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.75%  "
$ws.Range("D47").Value = "'2.203"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.05%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'3.538"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "'0.00000000336"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.76%  "
$ws.Range("D51").Value = "'0.07003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.96%  "
